$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.487.67'
$ws.Range('E2').Value = '  +4.24%  '

$ws.Range('D3').Value = '1.842.37'
$ws.Range('E3').Value = '  +3.75%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.032'
$ws.Range('E4').Value = '  +3.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.92'
$ws.Range('E5').Value = '  +4.58%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.026'
$ws.Range('E6').Value = '  +2.57%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4365'
$ws.Range('E7').Value = '  +3.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3725'
$ws.Range('E8').Value = '  +3.53%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07361'
$ws.Range('E9').Value = '  +3.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8739'
$ws.Range('E10').Value = '  +4.48%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.39'
$ws.Range('E11').Value = '  +4.78%  '

$ws.Range('D12').Value = '1.861.45'
$ws.Range('E12').Value = '  +4.68%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.474'
$ws.Range('E13').Value = '  +4.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.689'
$ws.Range('E14').Value = '  +3.68%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07147'
$ws.Range('E15').Value = '  +4.14%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.54'
$ws.Range('E16').Value = '  +4.49%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.030'
$ws.Range('E17').Value = '  +3.00%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009007'
$ws.Range('E18').Value = '  +4.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.024'
$ws.Range('E19').Value = '  +2.42%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.40'
$ws.Range('E20').Value = '  +3.40%  '

$ws.Range('D21').Value = '27.470.52'
$ws.Range('E21').Value = '  +4.13%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.230'
$ws.Range('E22').Value = '  +3.02%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.26'
$ws.Range('E23').Value = '  +3.17%  '

$ws.Range('D24').Value = '2.077.89'
$ws.Range('E24').Value = '  +3.91%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.74'
$ws.Range('E25').Value = '  +3.19%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.916'
$ws.Range('E26').Value = '  +6.67%  '

$ws.Range('E27').Value = '  +3.65%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.256'
$ws.Range('E28').Value = '  +3.66%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.922'
$ws.Range('E29').Value = '  +4.66%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.10'
$ws.Range('E30').Value = '  +1.31%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09058'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.204'
$ws.Range('E32').Value = '  +7.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7617'
$ws.Range('E33').Value = '  +4.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.481'
$ws.Range('E34').Value = '  +3.72%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.872'
$ws.Range('E35').Value = '  +5.22%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.027'
$ws.Range('E36').Value = '  +2.72%  '

$ws.Range('E37').Value = '  +5.28%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01965'
$ws.Range('E38').Value = '  +4.40%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05251'
$ws.Range('E39').Value = '  +2.45%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5162'
$ws.Range('E40').Value = '  +4.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.779'
$ws.Range('E41').Value = '  +6.63%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1663'
$ws.Range('E42').Value = '  +3.58%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.546'
$ws.Range('E43').Value = '  +3.57%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.486'
$ws.Range('E44').Value = '  +6.55%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.02'
$ws.Range('E45').Value = '  +4.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.63'
$ws.Range('E46').Value = '  +4.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.029'
$ws.Range('E47').Value = '  +2.93%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.694'
$ws.Range('E48').Value = '  +3.81%  '

$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4631'
$ws.Range('E49').Value = '  +4.10%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.899'
$ws.Range('E50').Value = '  +10.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06325'
$ws.Range('E51').Value = '  +2.45%  '
